$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old "Description of Core Business Processes" row (row 3).
# Rows below shift up: old row4->3, old row5->4, old row6->5.
$ws.Rows.Item(3).Delete()

# ---- Row 2: Description of Business/Project ----
$ws.Range("C2").Value = "An engaging and convincing story line is created describing the business, its challenges, the primary business problem to be solved and what resources are available, what constraints have to be considered, what costs and risks are involved. The description clearly follows the SMART steps."
$ws.Range("D2").Value = "All aspects of the business problem are outlined, connected and form an understandable description of the main issue(s) of the company/project and the type of solution that is being desired by the stakeholders. Key aspects like resources, constraints, costs, and other requirements provide a complete picture."
$ws.Range("E2").Value = "A basic outline of the business problem is provided. But important details are missing. Or the description is not concise but dives into too many irrelevant details are provided that distract from explaining the main business problem that needs to be solved."
$ws.Range("F2").Value = "The business problem is mentioned but not clearly outlined and the relevance of points made is not clear."
$ws.Range("G2").Value = "Missing or no work was submitted."
$ws.Range("C2:G2").Font.Color = 0
$ws.Range("C2:G2").VerticalAlignment = -4160
$ws.Range("C2:G2").WrapText = $true
$ws.Rows.Item(2).RowHeight = 119

# ---- Row 3: Description and Completeness of Possible Use Cases ----
$ws.Range("C3").Value = "Every use case is described clearly and concisely. Their relevance to solving the business problem clearly established. The outlines of their feasibility studies are sketched with timelines and a decision process to rank the use cases in terms feasibility and ability to create business value. "
$ws.Range("D3").Value = "A list of possible use cases is motivated well and connects directly with the business problem solution or aspects of it. A process of ranking the feasibility of each use case is outlined."
$ws.Range("E3").Value = "Use cases relevant to the business problem are mentioned. But their relative importance is not clear and/or key use cases are missing in the description."
$ws.Range("F3").Value = "The relevance of use cases described here in connection with the business problem at hand is not clear. "
$ws.Range("G3").Value = "Missing or no work was submitted."
$ws.Range("C3:G3").Font.Color = 0
$ws.Range("C3:G3").VerticalAlignment = -4160
$ws.Range("C3:G3").WrapText = $true
$ws.Rows.Item(3).RowHeight = 102

# ---- Row 4: Quality and Relevance of Visualizations ----
$ws.Range("C4").Value = "The visualization stands out in terms of quality, style, clarity, and its ability to focus the key messages in describing the business problem at hand. Color, graph choice, labeling, descriptions are thoughtfully and effectively "
$ws.Range("D4").Value = "The visualization used was relevant, to the point and added key information to illustrate the business problem succinctly and clearly. It used space, color, the choice of graphing style and other elements professionally and effectively to add to the overall message."
$ws.Range("E4").Value = "A visualization was added that was useful to communicate the overall description of the business problem but the visual could have been improved to make points more clearly or was not used strategically to emphasize key points in the presentation."
$ws.Range("F4").Value = "A visualization is included in the slide deck but only has limited relevance for the business problem discussion, key elements in the graph like axis labels are missing or hard to read."
$ws.Range("G4").Value = "Missing or no work was submitted."
$ws.Range("C4:G4").Font.Color = 0
$ws.Range("C4:G4").VerticalAlignment = -4160
$ws.Range("C4:G4").WrapText = $true
$ws.Rows.Item(4).RowHeight = 102

# Row 5 (was row 6, PowerPoint Presentation and Delivery) already holds its
# content/style and keeps height 119 - nothing further required there.

# ---- Page / view settings ----
$ws.PageSetup.Zoom = 49
$ws.PageSetup.Orientation = 2

$ws.Range("F4").Select()
